$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New case records to append starting at row 730 (columns A-I)
$rows = @(
    @("21TRD09386", "Bunner", "DUS UCM", "4510.111", "UCM", "No Contest", "Guilty", "$ 0", "$ 0"),
    @("21TRD09386", "Bunner", "TAIL LIGHTS-REAR LICENSE PLATE", "4513.05", "MM", "No Contest", "Guilty", "$ 0", "$ 0"),
    @("21TRD09386", "Bunner", "DUS UCM", "4510.111", "UCM", "No Contest", "Guilty", "$ 0", "$ 0"),
    @("21TRD09386", "Bunner", "TAIL LIGHTS-REAR LICENSE PLATE", "4513.05", "MM", "No Contest", "Guilty", "$ 0", "$ 0"),
    @("21TRD09386", "Bunner", "DUS UCM", "4510.111", "UCM", "No Contest", "Guilty", "$ 0", "$ 0"),
    @("21TRD09386", "Bunner", "TAIL LIGHTS-REAR LICENSE PLATE", "4513.05", "MM", "No Contest", "Guilty", "$ 0", "$ 0"),
    @("21TRC08418", "Bunner", "DRIVING IN MARKED LANES", "4511.33", "MM", "No Contest", "Guilty", "$ 0", "$ 0"),
    @("21TRC08418", "Bunner", "TURN AND STOP SIGNALS", "No Data", "MM", "No Contest", "Guilty", "$ 0", "$ 0"),
    @("21TRC08418", "Bunner", "OVI ALCOHOL / DRUGS 1ST", "4511.19A1A*", "M1", "No Contest", "Guilty", "$ 0", "$ 0"),
    @("22CRB00136", "Hemmeter", "DOMESTIC VIOLENCE", "2919.25(A)", "No Data", "Not Guilty"),
    @("22CRB00136", "Hemmeter", "ASSAULT - M1", "2903.13(A)", "No Data", "Not Guilty"),
    @("21CRB00626", "Hemmeter", "CRIMINAL MISCHIEF", "2909.07(A)(1)", "M3", "Not Guilty"),
    @("21CRB00626", "Hemmeter", "ASSAULT - M1", "2903.13(A)", "M1", "Not Guilty"),
    @("21CRB00626", "Hemmeter", "DISORDERLY CONDUCT", "2917.11A1", "MM", "Not Guilty"),
    @("22CRB00142", "Hemmeter", "THEFT / M1", "2913.02(A)(1)*", "M1", "Not Guilty")
)

$startRow = 730

# Cells whose literal text would otherwise be auto-converted by Excel into a
# number/currency value (e.g. "4510.111" -> 4510.111, "$ 0" -> currency 0)
# are collected here so they can be force-written as plain text afterwards,
# without leaving a quote-prefix / number-format style behind on the cell
# (the source workbook uses no per-cell styles at all).
$forceTextCells = @()

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $val = $vals[$c]
        if ($val -match '^[+-]?\d+(\.\d+)?$' -or $val -match '^\$\s*-?\d+(\.\d+)?$') {
            $forceTextCells += , @($r, ($c + 1), $val)
        } else {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
}

# Write each "numeric-looking" literal as a formula producing that exact
# string, then freeze it down to a literal value via copy/paste-special
# (values only) so no residual formula or special number format remains -
# just a plain text cell, matching how the rest of the sheet is encoded.
foreach ($entry in $forceTextCells) {
    $r = $entry[0]
    $c = $entry[1]
    $val = $entry[2]
    $quoted = $val.Replace('"', '""')
    $cell = $ws.Cells.Item($r, $c)
    $cell.Formula = '="' + $quoted + '"'
}

foreach ($entry in $forceTextCells) {
    $r = $entry[0]
    $c = $entry[1]
    $cell = $ws.Cells.Item($r, $c)
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
